# Updates the Price (D) and Volume(1h) (E) columns on the active sheet
# to reflect the latest scraped crypto values, per commit:
# "Updated symbol list on Sun Jan 22 03:55:39 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, new Price (column D, may be $null to skip), new Volume(1h) % (column E)
$updates = @(
    @{ Row = 2;  D = "300.49";     E = "-0.67%" },
    @{ Row = 3;  D = "38.10";      E = "8.76%" },
    @{ Row = 4;  D = "4.983";      E = "-3.04%" },
    @{ Row = 5;  D = "0.07719";    E = "-0.65%" },
    @{ Row = 6;  D = "2.193";      E = "-6.16%" },
    @{ Row = 7;  D = "7.957";      E = "-1.01%" },
    @{ Row = 8;  D = "3.995";      E = "1.17%" },
    @{ Row = 9;  D = "0.9156";     E = "-1.60%" },
    @{ Row = 10; D = "0.09095";    E = "-8.92%" },
    @{ Row = 11; D = "0.1793";     E = "0.23%" },
    @{ Row = 12; D = "0.08489";    E = "-1.55%" },
    @{ Row = 13; D = "0.03536";    E = "6.34%" },
    @{ Row = 14; D = "0.09935";    E = "0.19%" },
    @{ Row = 15; D = "0.001485";   E = "-1.03%" },
    @{ Row = 16; D = "0.005662";   E = "-1.43%" },
    @{ Row = 17; D = $null;        E = "0.42%" },
    @{ Row = 18; D = $null;        E = "4.02%" },
    @{ Row = 19; D = $null;        E = "2.88%" },
    @{ Row = 20; D = $null;        E = "-1.27%" },
    @{ Row = 21; D = "4.561";      E = "6.41%" },
    @{ Row = 22; D = $null;        E = "-2.90%" },
    @{ Row = 23; D = "0.04660";    E = "2.43%" },
    @{ Row = 24; D = $null;        E = "1.18%" },
    @{ Row = 25; D = "0.004435";   E = "1.27%" },
    @{ Row = 26; D = "0.0001302";  E = "0.06%" },
    @{ Row = 27; D = $null;        E = "40.12%" },
    @{ Row = 39; D = "0.01733";    E = "-3.16%" },
    @{ Row = 40; D = "0.04686";    E = "-2.26%" },
    @{ Row = 41; D = "0.007887";   E = "1.32%" },
    @{ Row = 42; D = $null;        E = "-1.77%" },
    @{ Row = 43; D = "0.007680";   E = "11.10%" },
    @{ Row = 44; D = "0.002304";   E = "9.58%" },
    @{ Row = 45; D = "0.009791";   E = "3.55%" },
    @{ Row = 46; D = "0.00006033"; E = "-1.30%" },
    @{ Row = 47; D = $null;        E = "0.07%" },
    @{ Row = 48; D = "8.621";      E = "184.81%" },
    @{ Row = 49; D = $null;        E = "34.82%" },
    @{ Row = 50; D = "0.00002103"; E = "0.07%" },
    @{ Row = 51; D = "0.0002003";  E = "0.07%" }
)

# Values must stay plain text (as in the source data), so force text
# formatting for the write and then restore the default "Normal" style so
# no stray style index is introduced on these previously-unstyled cells.
foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $cellD = $ws.Range("D$row")
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.Style = "Normal"
    }
    $cellE = $ws.Range("E$row")
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
    $cellE.Style = "Normal"
}
